# Refresh the "timestamp" column (Z) with the new run's timestamps.
# Every data row (2-112) in column Z gets overwritten; rows were written by
# the logger in consecutive batches that share the same timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$batches = @(
    @{ First = 2;   Last = 45;  Value = "2025-10-17T07:09:41.191455" },
    @{ First = 46;  Last = 55;  Value = "2025-10-17T07:09:41.248788" },
    @{ First = 56;  Last = 65;  Value = "2025-10-17T07:09:41.249785" },
    @{ First = 66;  Last = 74;  Value = "2025-10-17T07:09:41.250784" },
    @{ First = 75;  Last = 102; Value = "2025-10-17T07:09:41.307051" },
    @{ First = 103; Last = 112; Value = "2025-10-17T07:09:41.373676" }
)

foreach ($batch in $batches) {
    for ($row = $batch.First; $row -le $batch.Last; $row++) {
        $ws.Cells.Item($row, 26).Value = $batch.Value
    }
}
